$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.799.66"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "'3.479.67"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'415.95"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("D6").Value = "'129.65"
$ws.Range("E6").Value = "  +1.36%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.732"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("D11").Value = "'42.93"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").Value = "'9.58"
$ws.Range("E12").Value = "  +5.28%  "

$ws.Range("D13").Value = "'0.0000219"
$ws.Range("E13").Value = "  +8.22%  "

$ws.Range("D14").Value = "'4.013.82"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").Value = "'0.141"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "'20.62"
$ws.Range("E16").Value = "  -3.68%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.470.54"
$ws.Range("E17").Value = "  +1.97%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.87"
$ws.Range("E18").Value = "  +4.76%  "

$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "'62.794.54"
$ws.Range("E20").Value = "  +1.79%  "

$ws.Range("D21").Value = "'471.84"
$ws.Range("E21").Value = "  +5.84%  "

$ws.Range("D22").Value = "'91.18"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").Value = "'3.30"
$ws.Range("E23").Value = "  +3.51%  "

$ws.Range("D24").Value = "'13.43"
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("D25").Value = "'10.57"
$ws.Range("E25").Value = "  +22.37%  "

$ws.Range("D26").Value = "'3.32"
$ws.Range("E26").Value = "  +2.76%  "

$ws.Range("D27").Value = "'33.52"
$ws.Range("E27").Value = "  +2.20%  "

$ws.Range("D28").Value = "'4.81"
$ws.Range("E28").Value = "  +1.74%  "

$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'12.00"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.65"
$ws.Range("E31").Value = "  -2.43%  "

$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").Value = "'40.81"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").Value = "'58.23"
$ws.Range("E36").Value = "  +9.40%  "

$ws.Range("D37").Value = "'0.0492"
$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("E39").Value = "  +3.63%  "

$ws.Range("D40").Value = "'0.326"
$ws.Range("E40").Value = "  +3.02%  "

$ws.Range("D41").Value = "'3.37"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").Value = "'0.135"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("E43").Value = "  +5.55%  "

$ws.Range("D44").Value = "'145.36"
$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("D45").Value = "'2.09"
$ws.Range("E45").Value = "  +5.88%  "

$ws.Range("D46").Value = "'4.38"
$ws.Range("E46").Value = "  +4.84%  "

$ws.Range("D47").Value = "'2.40"
$ws.Range("E47").Value = "  +14.51%  "

$ws.Range("D48").Value = "'0.0₃0555"
$ws.Range("E48").Value = "  +36.72%  "

$ws.Range("D49").Value = "'16.48"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").Value = "'22.54"
$ws.Range("E50").Value = "  +1.38%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.140"
$ws.Range("E51").Value = "  +1.71%  "
